# Added RTC battery to BOM.  Fixed BOM quantities
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fixed BOM quantities
$ws.Range("F12").Value = 2
$ws.Range("F19").Value = 3

# Duplicate row 34 into a new row 35, preserving all formatting/styles
# (this is the RTC coin-cell battery row range with the dash/"-" source placeholder)
$ws.Rows(34).Copy()
$ws.Rows(35).Insert(-4121)

# Row 34 no longer needs the placeholder "-" note in column G
$ws.Range("G34").ClearContents()

# Row 35: fill in the new RTC battery part details
$ws.Range("C35").Value = "BATTERY LITHIUM 3V COIN 12.5MM"
$ws.Range("E35").Value = "P033-ND"
$ws.Range("G35").Clear()

# Restore the selection to match the saved view state
$ws.Range("G40").Select()
